$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J to make room for "param_frsmote__type"
$ws.Range("J1").EntireColumn.Insert()

# New header for J1
$ws.Range("J1").Value = 'param_frsmote__type'

# Row 2
$ws.Range("A2").Value = 31.11437217394511
$ws.Range("B2").Value = 0.1423787465209166
$ws.Range("C2").Value = 0.008661905924479166
$ws.Range("D2").Value = 0.002622076923139932
$ws.Range("E2").Value = 0.1
$ws.Range("F2").Value = 'lukasiewicz'
$ws.Range("G2").Value = 'linear'
$ws.Range("H2").Value = 'gaussian'
$ws.Range("I2").Value = 'minimum'
$ws.Range("J2").Value = 'itfrs'
$ws.Range("K2").Value = 'linear'
$ws.Range("L2").Value = 'minimum'
$ws.Range("M2").Value = 0.1
$ws.Range("N2").Value = 'rbf'
$ws.Range("O2").Value = '{''frsmote__gaussian_similarity_sigma'': 0.1, ''frsmote__lb_implicator_name'': ''lukasiewicz'', ''frsmote__lb_owa_method_name'': ''linear'', ''frsmote__similarity'': ''gaussian'', ''frsmote__similarity_tnorm'': ''minimum'', ''frsmote__type'': ''itfrs'', ''frsmote__ub_owa_method_name'': ''linear'', ''frsmote__ub_tnorm_name'': ''minimum'', ''svc__C'': 0.1, ''svc__kernel'': ''rbf''}'
$ws.Range("P2").Value = 0.6976744186046512
$ws.Range("Q2").Value = 0.6875
$ws.Range("R2").Value = 0.75
$ws.Range("S2").Value = 0.7117248062015503
$ws.Range("T2").Value = 0.0273815332382651
$ws.Range("U2").Value = 1

# Row 3
$ws.Range("A3").Value = 31.06540067990621
$ws.Range("B3").Value = 0.1872362441023085
$ws.Range("C3").Value = 0.007661819458007812
$ws.Range("D3").Value = 0.002053224730052947
$ws.Range("E3").Value = 0.1
$ws.Range("F3").Value = 'lukasiewicz'
$ws.Range("G3").Value = 'linear'
$ws.Range("H3").Value = 'gaussian'
$ws.Range("I3").Value = 'minimum'
$ws.Range("J3").Value = 'owafrs'
$ws.Range("K3").Value = 'linear'
$ws.Range("L3").Value = 'minimum'
$ws.Range("M3").Value = 0.1
$ws.Range("N3").Value = 'rbf'
$ws.Range("O3").Value = '{''frsmote__gaussian_similarity_sigma'': 0.1, ''frsmote__lb_implicator_name'': ''lukasiewicz'', ''frsmote__lb_owa_method_name'': ''linear'', ''frsmote__similarity'': ''gaussian'', ''frsmote__similarity_tnorm'': ''minimum'', ''frsmote__type'': ''owafrs'', ''frsmote__ub_owa_method_name'': ''linear'', ''frsmote__ub_tnorm_name'': ''minimum'', ''svc__C'': 0.1, ''svc__kernel'': ''rbf''}'
$ws.Range("P3").Value = 0.6976744186046512
$ws.Range("Q3").Value = 0.6875
$ws.Range("R3").Value = 0.75
$ws.Range("S3").Value = 0.7117248062015503
$ws.Range("T3").Value = 0.0273815332382651
$ws.Range("U3").Value = 1

# Row 4
$ws.Range("A4").Value = 30.97878360748291
$ws.Range("B4").Value = 0.1036403696230571
$ws.Range("C4").Value = 0.006995995839436849
$ws.Range("D4").Value = 0.002158913127057153
$ws.Range("E4").Value = 0.1
$ws.Range("F4").Value = 'lukasiewicz'
$ws.Range("G4").Value = 'linear'
$ws.Range("H4").Value = 'linear'
$ws.Range("I4").Value = 'minimum'
$ws.Range("J4").Value = 'itfrs'
$ws.Range("K4").Value = 'linear'
$ws.Range("L4").Value = 'minimum'
$ws.Range("M4").Value = 0.1
$ws.Range("N4").Value = 'rbf'
$ws.Range("O4").Value = '{''frsmote__gaussian_similarity_sigma'': 0.1, ''frsmote__lb_implicator_name'': ''lukasiewicz'', ''frsmote__lb_owa_method_name'': ''linear'', ''frsmote__similarity'': ''linear'', ''frsmote__similarity_tnorm'': ''minimum'', ''frsmote__type'': ''itfrs'', ''frsmote__ub_owa_method_name'': ''linear'', ''frsmote__ub_tnorm_name'': ''minimum'', ''svc__C'': 0.1, ''svc__kernel'': ''rbf''}'
$ws.Range("P4").Value = 0.6521739130434783
$ws.Range("Q4").Value = 0.6875
$ws.Range("R4").Value = 0.7368421052631579
$ws.Range("S4").Value = 0.692172006102212
$ws.Range("T4").Value = 0.03472315677086865
$ws.Range("U4").Value = 3

# Row 5
$ws.Range("A5").Value = 30.9464750289917
$ws.Range("B5").Value = 0.06694761131597743
$ws.Range("C5").Value = 0.005334854125976562
$ws.Range("D5").Value = 0.0004688699227739041
$ws.Range("E5").Value = 0.1
$ws.Range("F5").Value = 'lukasiewicz'
$ws.Range("G5").Value = 'linear'
$ws.Range("H5").Value = 'linear'
$ws.Range("I5").Value = 'minimum'
$ws.Range("J5").Value = 'owafrs'
$ws.Range("K5").Value = 'linear'
$ws.Range("L5").Value = 'minimum'
$ws.Range("M5").Value = 0.1
$ws.Range("N5").Value = 'rbf'
$ws.Range("O5").Value = '{''frsmote__gaussian_similarity_sigma'': 0.1, ''frsmote__lb_implicator_name'': ''lukasiewicz'', ''frsmote__lb_owa_method_name'': ''linear'', ''frsmote__similarity'': ''linear'', ''frsmote__similarity_tnorm'': ''minimum'', ''frsmote__type'': ''owafrs'', ''frsmote__ub_owa_method_name'': ''linear'', ''frsmote__ub_tnorm_name'': ''minimum'', ''svc__C'': 0.1, ''svc__kernel'': ''rbf''}'
$ws.Range("P5").Value = 0.6521739130434783
$ws.Range("Q5").Value = 0.6875
$ws.Range("R5").Value = 0.7368421052631579
$ws.Range("S5").Value = 0.692172006102212
$ws.Range("T5").Value = 0.03472315677086865
$ws.Range("U5").Value = 3
